$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is a daily price table: row 2 is always "today", and every
# subsequent row is the prior day. A new day's quote ("05-02-2026") is
# being published, so a fresh row is inserted at row 2 (pushing every
# existing row down by one) and the oldest row (old row 183) effectively
# reappears at the new row 184 with identical data.

# 1) Insert a blank row at row 2; rows 2..183 shift down to 3..184,
#    carrying their values/styles/number formats with them.
$ws.Rows("2:2").Insert()

# 2) Copy the cell formatting (styles/number formats) from row 3 (the
#    row that used to be row 2) into the new blank row 2 so D2 keeps the
#    "0.000" numeric style etc.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)

# 3) Populate row 2. Columns B, C, D, F repeat the former row 2's values
#    (same description/product code/price/link as the most recent
#    circular). Column A gets the new date. Column E keeps the former
#    row 2's circular date. Date-look-alike text is written with a
#    leading apostrophe so the COM layer stores it as literal text
#    instead of re-interpreting it as a date serial (matching how the
#    rest of the sheet stores dates as plain strings).
$ws.Range("B2").Value = $ws.Range("B3").Value2
$ws.Range("C2").Value = $ws.Range("C3").Value2
$ws.Range("D2").Value = $ws.Range("D3").Value2
$ws.Range("F2").Value = $ws.Range("F3").Value2

$ws.Range("A2").Value = "'05-02-2026"
$ws.Range("E2").Value = "'" + $ws.Range("E3").Value2

# 4) Row insertion does not carry the F-column hyperlink metadata past
#    the previous last row, so the newly revealed row 184 (old row 183)
#    needs its hyperlink re-attached, pointing at the same circular PDF
#    as its (identical) text already shows.
$ws.Hyperlinks.Add($ws.Cells.Item(184, 6), $ws.Range("F184").Value2)
